# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 8 (sheet data is sorted newest
# date first), pushing the existing rows 8 and 9 down to rows 9 and 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8:9 down to 9:10 by inserting a fresh blank row at 8.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the latest weekly observation.
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = 44524
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = "Otros"
$ws.Range("I8").Value = 100107002
$ws.Range("J8").Value = "Chirimoya"
$ws.Range("K8").Value = "Cultivar IV Región"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 23000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 23500
$ws.Range("Q8").Value = "$/caja 12 kilos"
$ws.Range("R8").Value = "Región de Coquimbo"
$ws.Range("S8").Value = 1958
$ws.Range("T8").Value = 12
